$d = $word.ActiveDocument

# The <id> element for this page previously held the "old" auto-generated
# id (p127v_a1) split across three runs: "<id>", "p127v_a1", "</id>".
# Replace it with the newly downloaded id "p127v_1", collapsed into the
# single "<id>p127v_1</id>" text run.
$d.Content.Find.Execute("<id>p127v_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p127v_1</id>", 2)
